# Organize data wrangling lessons
# - Split the single "Data wrangling" / "Data wrangling (more)" topics into
#   two more descriptive lessons.
# - Update the active cell selection to D8.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "Data wrangling: tidy data"
$ws.Range("D7").Value = "Data wrangling: relational data and factors"

$ws.Range("D8").Select()
